$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.655.89"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.597.56"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.50"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D12").Value = "1.821.99"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "1.626.28"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.00"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "26.641.39"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.69"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.98"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.32"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Value = "1.288.85"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.613"
$ws.Range("E36").Value = "  -8.74%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  +18.33%  "
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.44"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "1.734.36"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.72"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0508"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.38"
$ws.Range("E51").Value = "  -1.91%  "
